$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the "max_discount" column (X) data for every product row (2-485),
# matching the "light_schema"/max-discount value of 14 introduced by the
# "aggiunto light schema a foscarini" update. Row 1 already has the header.
$ws.Range("X2:X485").Value = 14
